# Applies the "Updated cryptos list" data refresh to the crypto price table.
# Column D (Price) and Column E (Volume(1h)) values are refreshed in place,
# and four coin rows (36-39, 45-46) are reordered/relabelled to reflect the
# new ranking returned by the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a literal text value into a cell without letting Excel
# --- reinterpret number-looking strings (e.g. "0.9993", "27.321.79") as
# --- actual numbers. We briefly mark the cell as Text ("@"), assign the
# --- literal string, then reset the cell style back to "Normal" so the
# --- cell keeps the plain/default style it had before (only its text
# --- content changes, same as the source workbook).
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# --- Column B / C: coin name + link swaps (plain text, never number-like) ---
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'

# --- Column D: Price (forced to literal text so exact formatting is kept) ---
Set-TextValue $ws.Range("D2") '27.321.79'
Set-TextValue $ws.Range("D3") '1.791.96'
Set-TextValue $ws.Range("D4") '0.9993'
Set-TextValue $ws.Range("D5") '338.47'
Set-TextValue $ws.Range("D6") '0.9981'
Set-TextValue $ws.Range("D7") '0.3937'
Set-TextValue $ws.Range("D8") '0.3455'
Set-TextValue $ws.Range("D9") '48.03'
Set-TextValue $ws.Range("D10") '1.199'
Set-TextValue $ws.Range("D11") '0.07480'
Set-TextValue $ws.Range("D12") '0.9964'
Set-TextValue $ws.Range("D13") '21.82'
Set-TextValue $ws.Range("D14") '6.474'
Set-TextValue $ws.Range("D15") '1.786.61'
Set-TextValue $ws.Range("D16") '7.112'
Set-TextValue $ws.Range("D17") '0.00001096'
Set-TextValue $ws.Range("D18") '0.06665'
Set-TextValue $ws.Range("D19") '84.05'
Set-TextValue $ws.Range("D20") '0.9991'
Set-TextValue $ws.Range("D21") '17.63'
Set-TextValue $ws.Range("D22") '6.527'
Set-TextValue $ws.Range("D23") '27.284.65'
Set-TextValue $ws.Range("D24") '12.37'
Set-TextValue $ws.Range("D25") '2.383'
Set-TextValue $ws.Range("D26") '1.508'
Set-TextValue $ws.Range("D27") '2.529'
Set-TextValue $ws.Range("D28") '21.21'
Set-TextValue $ws.Range("D30") '1.986.92'
Set-TextValue $ws.Range("D31") '134.73'
Set-TextValue $ws.Range("D32") '3.975'
Set-TextValue $ws.Range("D33") '6.047'
Set-TextValue $ws.Range("D34") '0.08715'
Set-TextValue $ws.Range("D35") '13.10'
Set-TextValue $ws.Range("D36") '5.441'
Set-TextValue $ws.Range("D37") '1.611'
Set-TextValue $ws.Range("D38") '0.6836'
Set-TextValue $ws.Range("D39") '0.02388'
Set-TextValue $ws.Range("D40") '0.06414'
Set-TextValue $ws.Range("D41") '0.2199'
Set-TextValue $ws.Range("D42") '1.248'
Set-TextValue $ws.Range("D43") '8.487'
Set-TextValue $ws.Range("D44") '14.30'
Set-TextValue $ws.Range("D45") '0.6426'
Set-TextValue $ws.Range("D46") '0.9969'
Set-TextValue $ws.Range("D47") '3.864'
Set-TextValue $ws.Range("D48") '2.138'
Set-TextValue $ws.Range("D49") '132.31'
Set-TextValue $ws.Range("D50") '0.07116'
Set-TextValue $ws.Range("D51") '78.95'

# --- Column E: Volume(1h) (already safe as text thanks to the padding
# --- spaces, but route through the same helper for consistency/safety) ---
Set-TextValue $ws.Range("E2") '  -0.60%  '
Set-TextValue $ws.Range("E3") '  -1.28%  '
Set-TextValue $ws.Range("E4") '  -0.45%  '
Set-TextValue $ws.Range("E5") '  -1.65%  '
Set-TextValue $ws.Range("E6") '  -0.29%  '
Set-TextValue $ws.Range("E7") '  +2.78%  '
Set-TextValue $ws.Range("E8") '  -2.67%  '
Set-TextValue $ws.Range("E9") '  -1.76%  '
Set-TextValue $ws.Range("E10") '  -2.95%  '
Set-TextValue $ws.Range("E11") '  -3.91%  '
Set-TextValue $ws.Range("E12") '  -0.56%  '
Set-TextValue $ws.Range("E13") '  -2.57%  '
Set-TextValue $ws.Range("E14") '  -2.01%  '
Set-TextValue $ws.Range("E15") '  -1.58%  '
Set-TextValue $ws.Range("E16") '  -1.52%  '
Set-TextValue $ws.Range("E17") '  -2.87%  '
Set-TextValue $ws.Range("E18") '  -0.97%  '
Set-TextValue $ws.Range("E19") '  -3.27%  '
Set-TextValue $ws.Range("E20") '  -0.18%  '
Set-TextValue $ws.Range("E21") '  -0.13%  '
Set-TextValue $ws.Range("E22") '  -0.35%  '
Set-TextValue $ws.Range("E23") '  -0.75%  '
Set-TextValue $ws.Range("E24") '  -5.78%  '
Set-TextValue $ws.Range("E25") '  -3.74%  '
Set-TextValue $ws.Range("E26") '  +2.83%  '
Set-TextValue $ws.Range("E27") '  -6.22%  '
Set-TextValue $ws.Range("E28") '  -4.38%  '
Set-TextValue $ws.Range("E29") '  +2.10%  '
Set-TextValue $ws.Range("E30") '  -1.59%  '
Set-TextValue $ws.Range("E31") '  -0.71%  '
Set-TextValue $ws.Range("E32") '  -2.13%  '
Set-TextValue $ws.Range("E33") '  -5.32%  '
Set-TextValue $ws.Range("E34") '  -1.15%  '
Set-TextValue $ws.Range("E35") '  -5.77%  '
Set-TextValue $ws.Range("E36") '  -3.48%  '
Set-TextValue $ws.Range("E37") '  -5.23%  '
Set-TextValue $ws.Range("E38") '  -2.53%  '
Set-TextValue $ws.Range("E39") '  -0.94%  '
Set-TextValue $ws.Range("E40") '  -1.47%  '
Set-TextValue $ws.Range("E42") '  -4.22%  '
Set-TextValue $ws.Range("E43") '  -6.39%  '
Set-TextValue $ws.Range("E44") '  -4.38%  '
Set-TextValue $ws.Range("E45") '  -3.05%  '
Set-TextValue $ws.Range("E46") '  -0.36%  '
Set-TextValue $ws.Range("E47") '  -2.67%  '
Set-TextValue $ws.Range("E48") '  -2.43%  '
Set-TextValue $ws.Range("E49") '  -0.74%  '
Set-TextValue $ws.Range("E50") '  -2.82%  '
Set-TextValue $ws.Range("E51") '  -2.43%  '

